$d = $word.ActiveDocument

# --- Step 1: remove the (hidden) _GoBack bookmark from its current spot.
#     It currently sits between the ". The number" run and the
#     " 1 can only divide..." run inside the NOTE- paragraph; deleting it
#     just drops the bookmarkStart/bookmarkEnd tags and leaves those two
#     runs adjacent, which is exactly what the target XML wants there.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: insert two new blank paragraphs (same formatting as the
#     existing trailing blank paragraph) right before that trailing
#     paragraph. Using InsertXML with a minimal package/paragraph (no
#     run) keeps the paragraph free of any stray runs, matching the
#     clean <w:p><w:pPr>...</w:pPr></w:p> shape in the target.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$blankParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + `
    '<w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' + `
    '</w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$insertionPoint.InsertXML($blankParaXml)
[void]$insertionPoint.InsertXML($blankParaXml)

# --- Step 3: re-add the _GoBack bookmark, collapsed at the very start of
#     what is now the final (formerly-trailing) blank paragraph.
#     A bookmark collapsed inside a run of trailing empty paragraphs can't
#     be placed directly (range resolution misbehaves right at the very
#     end of the document body), so temporarily type a placeholder
#     character into that paragraph, anchor the bookmark right before it,
#     then remove the placeholder again - the collapsed bookmark stays
#     put exactly where it was anchored.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.Text = "X"

$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmStart = $finalPara.Range.Start
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($bmStart, $bmStart + 1)
$placeholder.Delete()
